$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows for years 2000-2009 (original rows 2 through 11).
# This shifts the remaining rows (2010年.. 2018年) up to rows 2..7.
$ws.Range("A2:B11").EntireRow.Delete() | Out-Null

# Add the new row for 2021年 at the end (row 8).
$ws.Range("A8").Value = "2021年"
$ws.Range("B8").Value = 299

# Copy the style of an existing year cell (A2, which carries style index 1)
# onto the newly added A8 cell so its formatting matches the rest of the column.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$wb.Save()
